$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 54, shifting existing data (rows 54-175) down to (55-176)
$ws.Rows("54").Insert()

# Populate the new row 54 with the new price entry
$ws.Cells.Item(54, 1).Value = 10
$ws.Cells.Item(54, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(54, 3).Value = "La Araucanía"
$ws.Cells.Item(54, 4).Value = 44708
$ws.Cells.Item(54, 5).Value = 9
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100104
$ws.Cells.Item(54, 8).Value = "Frutos de pepita"
$ws.Cells.Item(54, 9).Value = 100104003
$ws.Cells.Item(54, 10).Value = "Membrillo"
$ws.Cells.Item(54, 11).Value = "Champion"
$ws.Cells.Item(54, 12).Value = "Primera"
$ws.Cells.Item(54, 13).Value = 170
$ws.Cells.Item(54, 14).Value = 10000
$ws.Cells.Item(54, 15).Value = 11000
$ws.Cells.Item(54, 16).Value = 10471
$ws.Cells.Item(54, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(54, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(54, 19).Value = 582
$ws.Cells.Item(54, 20).Value = 18
